# Elise runs through example code:
# swap the lat/lon columns (B <-> C), keeping per-cell formatting/styles and
# the column-width setting attached to the original "lat" column, then
# leave the selection on F3 (matching the state the workbook was saved in).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column at B. This shifts the original column B (lat values,
# incl. the bestFit column-width formatting and per-row styles) one column to
# the right, into C, and pushes the original column C (lon values) into D.
$ws.Range("B1:B9").Insert(-4161)

# Column B is now empty. Move the original lon data (now sitting in D, with
# its own per-row styles) into B.
$ws.Range("D1:D9").Cut($ws.Range("B1:B9"))

# Remove the now-empty scratch column D.
$ws.Range("D1:D9").Clear()

# Match the saved selection state.
$ws.Range("F3").Select()
